$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2236559139784946
$ws.Range("C2").Value = 0.5053763440860215
$ws.Range("J2").Value = 0.01505376344086022
$ws.Range("P2").Value = 0.1870967741935484
$ws.Range("S2").Value = 0.06881720430107527

# Row 3
$ws.Range("B3").Value = 0.007936507936507936
$ws.Range("C3").Value = 0.03571428571428571
$ws.Range("J3").Value = 0.02380952380952381
$ws.Range("P3").Value = 0.7619047619047619
$ws.Range("S3").Value = 0.1706349206349206

# Row 4
$ws.Range("J4").Value = 0.1388888888888889
$ws.Range("P4").Value = 0.6388888888888888
$ws.Range("S4").Value = 0.2222222222222222

# Row 5
$ws.Range("J5").Value = 0.5
$ws.Range("P5").Value = 0.5

# Row 6
$ws.Range("B6").Value = 0.07624633431085044
$ws.Range("D6").Value = 0.02052785923753666
$ws.Range("F6").Value = 0.07038123167155426
$ws.Range("J6").Value = 0.281524926686217
$ws.Range("O6").Value = 0.03519061583577713
$ws.Range("Q6").Value = 0.1642228739002932
$ws.Range("R6").Value = 0.05571847507331378
$ws.Range("S6").Value = 0.2961876832844575

# Row 7
$ws.Range("B7").Value = 0.1270491803278689
$ws.Range("D7").Value = 0.01639344262295082
$ws.Range("F7").Value = 0.04918032786885246
$ws.Range("J7").Value = 0.139344262295082
$ws.Range("O7").Value = 0.04098360655737705
$ws.Range("Q7").Value = 0.1516393442622951
$ws.Range("R7").Value = 0.1188524590163934
$ws.Range("S7").Value = 0.3565573770491803

# Row 8
$ws.Range("B8").Value = 0.1104868913857678
$ws.Range("D8").Value = 0.02434456928838951
$ws.Range("F8").Value = 0.06554307116104868
$ws.Range("J8").Value = 0.1179775280898876
$ws.Range("O8").Value = 0.02621722846441948
$ws.Range("Q8").Value = 0.2209737827715356
$ws.Range("R8").Value = 0.1161048689138577
$ws.Range("S8").Value = 0.3183520599250936

# Row 9
$ws.Range("B9").Value = 0.09266409266409266
$ws.Range("D9").Value = 0.04633204633204633
$ws.Range("F9").Value = 0.07335907335907337
$ws.Range("J9").Value = 0.1196911196911197
$ws.Range("O9").Value = 0.03861003861003861
$ws.Range("Q9").Value = 0.1814671814671815
$ws.Range("R9").Value = 0.1003861003861004
$ws.Range("S9").Value = 0.3474903474903475

# Row 10
$ws.Range("B10").Value = 0.1066195048004042
$ws.Range("D10").Value = 0.01970692268822638
$ws.Range("E10").Value = 0.001515917129863567
$ws.Range("F10").Value = 0.06114199090449722
$ws.Range("J10").Value = 0.1217786760990399
$ws.Range("O10").Value = 0.02425467407781708
$ws.Range("Q10").Value = 0.2531581606872157
$ws.Range("R10").Value = 0.08943911066195048
$ws.Range("S10").Value = 0.3223850429509854

# Row 11
$ws.Range("F11").Value = 0.002398081534772182
$ws.Range("G11").Value = 0.1870503597122302
$ws.Range("J11").Value = 0.1247002398081535
$ws.Range("K11").Value = 0.2350119904076739
$ws.Range("L11").Value = 0.4388489208633093
$ws.Range("S11").Value = 0.01199040767386091

# Row 12
$ws.Range("G12").Value = 0.6878306878306878
$ws.Range("J12").Value = 0.2592592592592592
$ws.Range("K12").Value = 0.01587301587301587
$ws.Range("L12").Value = 0.02645502645502645
$ws.Range("S12").Value = 0.01058201058201058

# Row 13
$ws.Range("G13").Value = 0.5714285714285714
$ws.Range("J13").Value = 0.3857142857142857
$ws.Range("S13").Value = 0.04285714285714286

# Row 15
$ws.Range("F15").Value = 0.02380952380952381
$ws.Range("H15").Value = 0.1507936507936508
$ws.Range("I15").Value = 0.03968253968253968
$ws.Range("J15").Value = 0.3941798941798942
$ws.Range("K15").Value = 0.07407407407407407
$ws.Range("M15").Value = 0.005291005291005291
$ws.Range("O15").Value = 0.08994708994708994
$ws.Range("S15").Value = 0.2222222222222222

# Row 16
$ws.Range("F16").Value = 0.0255591054313099
$ws.Range("H16").Value = 0.1277955271565495
$ws.Range("I16").Value = 0.08626198083067092
$ws.Range("J16").Value = 0.4249201277955272
$ws.Range("K16").Value = 0.1054313099041534
$ws.Range("M16").Value = 0.04472843450479233
$ws.Range("O16").Value = 0.08945686900958466
$ws.Range("S16").Value = 0.09584664536741214

# Row 17
$ws.Range("F17").Value = 0.03078982597054886
$ws.Range("H17").Value = 0.1686746987951807
$ws.Range("I17").Value = 0.108433734939759
$ws.Range("J17").Value = 0.461847389558233
$ws.Range("K17").Value = 0.08165997322623829
$ws.Range("M17").Value = 0.01338688085676037
$ws.Range("N17").Value = 0.002677376171352075
$ws.Range("O17").Value = 0.05890227576974565
$ws.Range("S17").Value = 0.07362784471218206

# Row 18
$ws.Range("F18").Value = 0.02922077922077922
$ws.Range("H18").Value = 0.1331168831168831
$ws.Range("I18").Value = 0.1168831168831169
$ws.Range("J18").Value = 0.4577922077922078
$ws.Range("K18").Value = 0.09415584415584416
$ws.Range("M18").Value = 0.02597402597402598
$ws.Range("O18").Value = 0.07467532467532467
$ws.Range("S18").Value = 0.06818181818181818

# Row 19
$ws.Range("F19").Value = 0.026512576478586
$ws.Range("H19").Value = 0.1876274643099932
$ws.Range("I19").Value = 0.06866077498300475
$ws.Range("J19").Value = 0.4194425560842964
$ws.Range("K19").Value = 0.1046906866077498
$ws.Range("M19").Value = 0.026512576478586
$ws.Range("N19").Value = 0.001359619306594154
$ws.Range("O19").Value = 0.08565601631543168
$ws.Range("S19").Value = 0.07953772943575799
